$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.304.72'
$ws.Range("E2").Value = '  -2.63%  '
$ws.Range("D3").Value = '1.852.98'
$ws.Range("E3").Value = '  -3.26%  '
$ws.Range("D5").Value = "'326.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = "'0.4554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.97%  '
$ws.Range("D8").Value = "'0.3886"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.34%  '
$ws.Range("D9").Value = "'48.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -9.17%  '
$ws.Range("D10").Value = "'0.07929"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.39%  '
$ws.Range("D11").Value = "'1.014"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.80%  '
$ws.Range("D12").Value = "'21.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").Value = '1.866.19'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = "'5.910"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = "'7.155"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = "'86.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.00%  '
$ws.Range("D18").Value = "'0.06594"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = "'0.00001028"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").Value = "'17.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.11%  '
$ws.Range("D22").Value = "'5.496"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.57%  '
$ws.Range("D23").Value = '27.311.47'
$ws.Range("E23").Value = '  -2.63%  '
$ws.Range("D24").Value = "'10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.15%  '
$ws.Range("D25").Value = "'2.293"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").Value = '2.085.18'
$ws.Range("E26").Value = '  -0.92%  '
$ws.Range("D27").Value = "'154.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").Value = "'19.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").Value = "'2.062"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.17%  '
$ws.Range("D30").Value = "'5.461"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.14%  '
$ws.Range("D31").Value = "'121.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'0.9466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = "'0.09350"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("D34").Value = "'1.442"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").Value = "'3.590"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("E36").Value = '  -4.81%  '
$ws.Range("D37").Value = "'0.06024"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").Value = "'0.02229"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.18%  '
$ws.Range("D39").Value = "'1.211"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("D40").Value = "'8.056"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.05%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = "'0.5927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("D43").Value = "'0.1886"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").Value = "'10.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.66%  '
$ws.Range("D45").Value = "'1.283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.81%  '
$ws.Range("D46").Value = "'0.5606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("D47").Value = "'12.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.13%  '
$ws.Range("D48").Value = "'3.384"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("D49").Value = "'1.912"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.17%  '
$ws.Range("D50").Value = "'0.06736"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = "'108.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.21%  '
